# adding revisions to folder structure
# - Add a new note "STARTING WITH THIS DATASET. " in cell D2 (new shared string)
# - Select D2 (matches the saved cursor/selection position)
# - Widen column A slightly to better fit the existing content

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New cell content for D2
$ws.Range("D2").Value = "STARTING WITH THIS DATASET. "

# Make D2 the active selection, as recorded when the file was saved
[void]$ws.Range("D2").Select()

# Column A grew wider in the saved file
$ws.Columns.Item(1).ColumnWidth = 46.33
